$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.16690719127655
$ws.Range("B1").Value = 4.087239742279053
$ws.Range("C1").Value = 3.422869443893433
$ws.Range("D1").Value = 2.406485557556152
$ws.Range("E1").Value = 2.171691417694092
